$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 356; existing rows 356-392 shift down to 357-393.
$ws.Rows.Item(356).Insert()

# Populate the newly inserted row 356 with the new record's data.
$ws.Cells.Item(356, 1).Value = 5
$ws.Cells.Item(356, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(356, 3).Value = "Maule"
$ws.Cells.Item(356, 4).Value = 44858
$ws.Cells.Item(356, 5).Value = 7
$ws.Cells.Item(356, 6).Value = 100112006
$ws.Cells.Item(356, 7).Value = "Repollo"
$ws.Cells.Item(356, 8).Value = "Crespo record"
$ws.Cells.Item(356, 9).Value = "Primera"
$ws.Cells.Item(356, 10).Value = 3000
$ws.Cells.Item(356, 11).Value = 1800
$ws.Cells.Item(356, 12).Value = 1800
$ws.Cells.Item(356, 13).Value = 1800
$ws.Cells.Item(356, 14).Value = "$/unidad"
$ws.Cells.Item(356, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(356, 16).Value = 1800
$ws.Cells.Item(356, 17).Value = 1
$ws.Cells.Item(356, 18).Value = "Hortaliza"
